$wb = $excel.ActiveWorkbook

# --- ACCOUNT sheet: delete the ORGID row (row 4) ---
$acc = $wb.Worksheets.Item("ACCOUNT")
$acc.Activate()
$acc.Rows.Item(4).Select()

# Grab the two existing comments (anchored at E10 and E12) before they get
# re-parented so we can re-create them one row higher, matching the content
# that shifts up when the row is removed.
$comment1 = $acc.Range("E10").Comment
$comment1Text = $comment1.Text()
$comment2 = $acc.Range("E12").Comment
$comment2Text = $comment2.Text()

$acc.Rows.Item(4).Delete()

$comment1.Delete()
$acc.Range("E9").AddComment($comment1Text)

$comment2.Delete()
$acc.Range("E11").AddComment($comment2Text)

# --- SCHEDULE sheet: view no longer frozen/scrolled to A3, selection moved to B14 ---
$sched = $wb.Worksheets.Item("SCHEDULE")
$sched.Activate()
$sched.Range("B14").Select()

# --- STATISTIC sheet: no longer the tab that's selected when file is saved ---
$stat = $wb.Worksheets.Item("STATISTIC")
$stat.Activate()

# --- Re-activate ACCOUNT so it is the tab/sheet shown when the workbook opens ---
$acc.Activate()
